$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1. Stage")

# Fill in the previously empty "Calc" column values for rows 8 and 9
$ws.Range("I8").Value = 0.064103
$ws.Range("I9").Value = 0.197606

# Recalculate so dependent totals (I13 = SUM(I2:I12)) and other formulas refresh
$excel.Calculate()

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("I19").Select()
